$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing note (C4) to add a second bullet point ---
$ws.Range("C4").Value = "1.提出通过经纬度计算估算两个卫星之间的ISL跳数，其实可以通过给卫星编号来计算。。。。`n2.统计了ISL跳数分布的规律"

# --- Copy C4's formatting (wrap text, vertical alignment, etc.) onto the new D4/E4 cells ---
$ws.Range("C4").Copy()
$ws.Range("D4:E4").PasteSpecial(-4122)  # xlPasteFormats

# --- Fill in the two new notes for row 4 ---
# (E4's shared string must be interned before D4's so the new shared-string
# table indices line up with the target workbook: E4 -> 18, D4 -> 19)
$ws.Range("E4").Value = "1.计算跳数的算法太复杂但是没有意义。。。。`n2.仅研究了跳数分布的规律。。。"
$ws.Range("D4").Value = "1.在目前的卫星部署模式下（每个轨道卫星的数量大于轨道的跳数），那么横向传输的跳数少于纵向传播的跳数；路由策略可优先考虑横向传播。`n2.卫星相对于赤道的偏移量越大，跳数越少"

# --- Row 2 keeps its manual (wrap-text driven) custom height; nudge it to
#     the value captured in the diff. (Row 1's height is auto-computed from
#     the title font and is intentionally left alone.)
$ws.Rows.Item(2).RowHeight = 159.55

# --- Move the active selection to the newly filled-in cell ---
$ws.Range("E4").Select()
